$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(310.0, "Sunday, Jan 08", "7:45 PM", "LO6193", "Phuket", "(HKT)", "LOT ", "B788", "(SP-LRD)", "7:56 PM", $null, "0 hours, 11 minutes", $null)
    ,@(311.0, "Sunday, Jan 08", "7:50 PM", "LX1353", "Zurich", "(ZRH)", "Air Baltic ", "BCS3", "(YL-ABF)", "8:16 PM", $null, "0 hours, 26 minutes", $null)
    ,@(312.0, "Sunday, Jan 08", "7:55 PM", "LO137", "Istanbul", "(IST)", "LOT (Grzeski Livery) ", "E195", "(SP-LNB)", "8:05 PM", $null, "0 hours, 10 minutes", $null)
    ,@(313.0, "Sunday, Jan 08", "7:55 PM", "LO3825", "Gdansk", "(GDN)", "LOT ", "E75S", "(SP-LIQ)", "8:08 PM", $null, "0 hours, 13 minutes", $null)
    ,@(314.0, "Sunday, Jan 08", "7:55 PM", "LO3859", "Wroclaw", "(WRO)", "LOT ", "E75S", "(SP-LIC)", "8:03 PM", $null, "0 hours, 8 minutes", $null)
    ,@(315.0, "Sunday, Jan 08", "7:55 PM", "LO3921", "Krakow", "(KRK)", "LOT ", "E170", "(SP-LDF)", "8:10 PM", $null, "0 hours, 15 minutes", $null)
    ,@(316.0, "Sunday, Jan 08", "8:00 PM", "LO3985", "Zielona Gora", "(IEG)", "LOT (Star Alliance Livery) ", "E75S", "(SP-LIO)", "8:07 PM", $null, "0 hours, 7 minutes", $null)
    ,@(317.0, "Sunday, Jan 08", "8:15 PM", "LO231", "Brussels", "(BRU)", "LOT ", "E195", "(SP-LNO)", "8:33 PM", $null, "0 hours, 18 minutes", $null)
    ,@(318.0, "Sunday, Jan 08", "8:15 PM", "LO333", "Paris", "(CDG)", "LOT (Retro Livery) ", "E75S", "(SP-LIM)", "8:24 PM", $null, "0 hours, 9 minutes", $null)
    ,@(319.0, "Sunday, Jan 08", "8:15 PM", "LO285", "London", "(LHR)", "LOT ", "B38M", "(SP-LVA)", "9:30 PM", $null, "1 hours, 15 minutes", $null)
    ,@(320.0, "Sunday, Jan 08", "8:25 PM", "LO405", "Dusseldorf", "(DUS)", "LOT ", "E190", "(SP-LMF)", "8:38 PM", $null, "0 hours, 13 minutes", $null)
    ,@(321.0, "Sunday, Jan 08", "8:30 PM", "LO529", "Prague", "(PRG)", "LOT ", "E75S", "(SP-LIN)", "8:28 PM", $null, "0 hours, -2 minutes", $null)
    ,@(322.0, "Sunday, Jan 08", "8:35 PM", "LO401", "Hamburg", "(HAM)", "LOT ", "E170", "(SP-LDH)", "8:43 PM", $null, "0 hours, 8 minutes", $null)
    ,@(323.0, "Sunday, Jan 08", "8:35 PM", "LO455", "Stockholm", "(ARN)", "LOT ", "E190", "(SP-LMG)", "8:52 PM", $null, "0 hours, 17 minutes", $null)
    ,@(324.0, "Sunday, Jan 08", "8:35 PM", "LO495", "Gothenburg", "(GOT)", "LOT ", "E195", "(SP-LND)", "8:57 PM", $null, "0 hours, 22 minutes", $null)
    ,@(325.0, "Sunday, Jan 08", "8:35 PM", "LO719", "Baku", "(GYD)", "LOT ", "B38M", "(SP-LVB)", "9:00 PM", $null, "0 hours, 25 minutes", $null)
    ,@(326.0, "Sunday, Jan 08", "8:40 PM", "LO383", "Frankfurt", "(FRA)", "LOT ", "E75S", "(SP-LIA)", "8:54 PM", $null, "0 hours, 14 minutes", $null)
    ,@(327.0, "Sunday, Jan 08", "8:40 PM", "LO459", "Copenhagen", "(CPH)", "LOT ", "E195", "(SP-LNG)", "8:50 PM", $null, "0 hours, 10 minutes", $null)
    ,@(328.0, "Sunday, Jan 08", "9:00 PM", "W61539", "Reykjavik", "(KEF)", "Wizz Air ", "A21N", "(HA-LZE)", "9:27 PM", $null, "0 hours, 27 minutes", $null)
    ,@(329.0, "Sunday, Jan 08", "9:15 PM", "FZ1830", "Dubai", "(DXB)", "flydubai ", "B38M", "(A6-FMM)", "10:11 PM", $null, "0 hours, 56 minutes", $null)
    ,@(330.0, "Sunday, Jan 08", "9:50 PM", "LO727", "Yerevan", "(EVN)", "LOT ", "B738", "(SP-LWA)", "10:05 PM", $null, "0 hours, 15 minutes", $null)
)

$startRow = 311
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $val = $data[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

Write-Output "Added $($rows.Count) rows starting at row $startRow"
